# "Generate Report for Archive"
# Update the localization status of the two tracked files from
# "Ready for handoff" to "In Translation" on every sheet that surfaces
# it (Overview summary + the per-locale detail sheets), then refresh the
# affected status-column widths to their new best-fit size.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# --- zh-cn detail sheet: Status column (C) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

# --- de-de detail sheet: Status column (C) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Resize the status columns to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
